# Applies the "Updated symbol list" commit: refreshes crypto price/volume
# data for several rows, and shifts coin rows 10-18 down by one (inserting
# "One" at row 10 and dropping the old row-18 duplicate of "One").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "265.08"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "6.229"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06173"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.560"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.735"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.357"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8169"
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.01354"
$ws.Range("E10").Value = "9OneONE"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1608"
$ws.Range("E11").Value = "10WazirXWRX"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08219"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03383"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03163"
$ws.Range("E14").Value = "13BitrueCoinBTR"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09257"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.894"
$ws.Range("E16").Value = "15MCDexMCB"
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001702"
$ws.Range("E17").Value = "16BitForexTokenBF"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.04853"
$ws.Range("E18").Value = "17CoinExTokenCET"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006237"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.006273"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.001099"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.698"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.263"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0002681"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04592"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007250"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003225"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01044"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006118"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.7699"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1963"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002100"
